# Adapt column header formatting to respective input file names (#7)
#  - headers that used to carry an "_old" / "_new" suffix now carry the
#    file-format-version suffix they actually belong to ("_FV2210" / "_FV2304")
#  - wrap the data range in a proper Excel Table ("Table1")
#  - freeze the header row so it stays visible while scrolling

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the column headers in row 1 --------------------------------
# Columns A-J ("...-_old") -> "..._FV2210"
# Column K ("diff") is left untouched
# Columns L-U ("..._new") -> "..._FV2304"

for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2210")
}

for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2304")
}

# --- 2. Turn the used range into a proper Table ----------------------------

$dataRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------

$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
